# cryptos.xlsx - "Updated cryptos list ... with GitHub Actions"
# Refresh the Price (column D) and Volume(1h) (column E) columns with the
# latest scraped quote for every coin row (2-51). Both columns hold plain
# text in the workbook (not numbers/percentages), so any value that would
# otherwise be auto-recognised by Excel as a number (e.g. "1.003") is
# written to a cell pre-formatted as Text ("@") to keep it a literal string,
# matching the original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.760.85"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "1.699.97"

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.69"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3931"
$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.503"
$ws.Range("E9").Value = "  -3.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "54.03"
$ws.Range("E10").Value = "  -2.09%  "

$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08886"
$ws.Range("E12").Value = "  +0.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.241"
$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.37"
$ws.Range("E14").Value = "  -0.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.038"
$ws.Range("E15").Value = "  +5.45%  "

$ws.Range("E16").Value = "  -0.51%  "

$ws.Range("D17").Value = "1.702.11"
$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "100.19"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07006"
$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.62"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.030"
$ws.Range("E21").Value = "  +1.40%  "

$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.45"
$ws.Range("E23").Value = "  +1.89%  "

$ws.Range("D24").Value = "24.745.87"
$ws.Range("E24").Value = "  +0.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.251"
$ws.Range("E25").Value = "  +9.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.354"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.77"
$ws.Range("E27").Value = "  +1.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.99"
$ws.Range("E28").Value = "  +0.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "136.23"
$ws.Range("E29").Value = "  +1.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.170"
$ws.Range("E30").Value = "  -1.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.752"
$ws.Range("E31").Value = "  +1.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08734"
$ws.Range("E32").Value = "  +1.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.073"
$ws.Range("E33").Value = "  -4.03%  "

$ws.Range("E34").Value = "  -3.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.23"
$ws.Range("E35").Value = "  +0.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.958"
$ws.Range("E36").Value = "  -1.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2745"
$ws.Range("E37").Value = "  -0.92%  "

$ws.Range("E38").Value = "  -3.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09191"
$ws.Range("E39").Value = "  +1.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02733"
$ws.Range("E40").Value = "  -1.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.461"
$ws.Range("E41").Value = "  -0.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7677"
$ws.Range("E42").Value = "  -1.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.87"
$ws.Range("E43").Value = "  +2.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7173"
$ws.Range("E44").Value = "  -1.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.570"
$ws.Range("E45").Value = "  +1.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.219"
$ws.Range("E46").Value = "  +0.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "140.86"
$ws.Range("E48").Value = "  -0.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.308"
$ws.Range("E49").Value = "  +0.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07979"
$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "90.50"
$ws.Range("E51").Value = "  +2.49%  "
